$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 23 already holds the "source" formatting (date-column style, text-column
# style, and the wrap/vertical-center style used for the Achievements /
# Reflection columns). Copy that formatting down into the three new diary
# rows (24-26) before filling in their values, so the new rows look like the
# existing ones.
$ws.Range("A23:G23").Copy() | Out-Null
$ws.Range("A24:G26").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Row 24 - 2020-02-19
$ws.Cells.Item(24, 1).Value = 43880
$ws.Cells.Item(24, 2).Value = "2:30PM - 5:00 PM"
$ws.Cells.Item(24, 3).Value = "Chris Zhang, Nicolas Grantham, and Hyun Jay Yang"
$ws.Cells.Item(24, 4).Value = "Finish and resubmit assignment 2"
$ws.Cells.Item(24, 5).Value = "We rewrote the document using higher level abstractions and went our discovery process more throughly"
$ws.Cells.Item(24, 6).Value = "Structuring the document into paragraphs and making links to our diagrams makes our document more understandable"
$ws.Cells.Item(24, 7).Value = "Feeling good overall"
$ws.Rows(24).RowHeight = 63

# Row 26 - 2020-02-21 (entered before row 25)
$ws.Cells.Item(26, 1).Value = 43882
$ws.Cells.Item(26, 2).Value = "2:00PM - 7:10 PM"
$ws.Cells.Item(26, 3).Value = "Chris Zhang, Nicolas Grantham, and Hyun Jay Yang"
$ws.Cells.Item(26, 7).Value = "Exhausted"
$ws.Cells.Item(26, 5).Value = "Explained the social context, identified interesting pull requests and issues, and explained the architecture of our project in a concise document"
$ws.Cells.Item(26, 4).Value = "Finish and deliver our project assignment #4"
$ws.Cells.Item(26, 6).Value = "Since we had already worked on the essential features, we already had a general understanding of our project, so it was easier to understand the architecture because we knew the exact routes we had to study"
$ws.Rows(26).RowHeight = 110.25

# Row 25 - 2020-02-20
$ws.Cells.Item(25, 1).Value = 43881
$ws.Cells.Item(25, 2).Value = "5:00PM - 7:00 PM"
$ws.Cells.Item(25, 3).Value = "N/A"
$ws.Cells.Item(25, 4).Value = "Learn new expert key practices, what is social context, and how does architecture can help understanding code"
$ws.Cells.Item(25, 5).Value = "Understood what social context is and how it might affect our decisions when choosing a project, and  how professionals use architecture as a comprehension tool"
$ws.Cells.Item(25, 7).Value = "Feeling ok"
$ws.Cells.Item(25, 6).Value = "It was good to hear our guest speakers give suggestions on how to introduce new members to the project and team, and how they guide them in the proper way to contribute"
$ws.Rows(25).RowHeight = 94.5

# Match the author's on-screen position when they saved: scrolled down so row
# 22 is at the top, with cell A24 selected.
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 22
$win.ScrollColumn = 1
$ws.Range("A24").Select() | Out-Null
